# "Invalid test case 2"
#
# - Adds a third worksheet "InvalidLogin" (after "ValidLogin") that mirrors
#   the existing "ValidLogin" sheet's username/password layout, but with
#   invalid credentials ("abc" / "xyz") used to exercise a negative login
#   test case.
# - The new sheet becomes the active tab, and its column A is autofit to
#   its content (matches the commit's width/bestFit metadata on sheet3).
# - The previously active "ValidLogin" sheet's selection is updated to
#   A1:B2 (no longer the lone active cell "B3"), and it stops being the
#   active/selected tab now that "InvalidLogin" is active.

$wb = $excel.ActiveWorkbook

# --- Update the "ValidLogin" sheet's selection -----------------------------
$wsValid = $wb.Worksheets.Item("ValidLogin")
$wsValid.Range("A1:B2").Select()

# --- Add the new "InvalidLogin" worksheet at the end of the tab strip -----
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$wsInvalid  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsInvalid.Name = "InvalidLogin"

# --- Populate it with header + invalid credential values -------------------
$wsInvalid.Range("A1").Value = "username"
$wsInvalid.Range("B1").Value = "password"
$wsInvalid.Range("A2").Value = "abc"
$wsInvalid.Range("B2").Value = "xyz"

# --- Autofit column A (matches the recorded bestFit column width) ---------
$wsInvalid.Columns.Item(1).AutoFit() | Out-Null

# --- Leave selection on B4, the cell just below the data, as recorded -----
$wsInvalid.Range("B4").Select()
